# Updates "F" column (attendance / view counts) figures across the four
# worksheets of the 杭州-漫展信息 workbook, matching the data refresh
# recorded in the commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Map of sheet name -> { row number -> new F-column value }
$updates = @{
    "展览" = @{
        2  = 833
        4  = 1711
        5  = 14
        6  = 522
        7  = 2142
        8  = 1342
        9  = 2023
        10 = 957
        13 = 636
        15 = 3785
        17 = 357
        18 = 2844
        19 = 743
        20 = 137
        22 = 63
        23 = 1984
        24 = 1157
        25 = 1751
        26 = 362
        27 = 191
        28 = 7950
        29 = 5459
        30 = 342
        31 = 161
        32 = 738
        33 = 750
        34 = 3455
        36 = 937
        37 = 369
        38 = 172
        39 = 146
        40 = 4571
        41 = 796
        42 = 49
        43 = 372
    }
    "演出" = @{
        15 = 104
        17 = 113
        18 = 144
        25 = 26
    }
    "本地生活" = @{
        2 = 8155
        3 = 349
        4 = 1209
    }
    "全部类型" = @{
        2  = 8155
        3  = 833
        4  = 349
        5  = 1209
        8  = 1711
        9  = 522
        10 = 1342
        11 = 957
        14 = 3785
        15 = 357
        16 = 2844
        17 = 743
        18 = 137
        20 = 1984
        25 = 1157
        27 = 1751
        28 = 104
        29 = 362
        30 = 191
        31 = 7950
        32 = 5459
        34 = 342
        35 = 738
        36 = 750
        37 = 3455
        39 = 937
        40 = 369
        41 = 172
        43 = 146
        44 = 4571
        45 = 796
        46 = 49
        47 = 372
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowMap = $updates[$sheetName]
    foreach ($row in $rowMap.Keys) {
        $ws.Cells.Item($row, 6).Value = $rowMap[$row]
    }
}
